$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 242.83333
$ws.Range("I9").Value = 161.5
$ws.Range("K9").Value = 161.5
$ws.Range("M9").Value = 7.5
$ws.Range("H19").Value = 2561.8235
$ws.Range("I19").Value = 2349.9167
$ws.Range("K19").Value = 2349.9167
$ws.Range("M19").Value = -2174.9167
$ws.Range("H107").Value = 475.25
$ws.Range("I107").Value = 426.83334
$ws.Range("J107").Value = 620.5
$ws.Range("K107").Value = 426.83334
$ws.Range("L107").Value = 620.5
$ws.Range("M107").Value = 1493.16666
$ws.Range("N107").Value = -4460.5
$ws.Range("H137").Value = 1252233.4
$ws.Range("I137").Value = 3824.6667
$ws.Range("J137").Value = 2322298
$ws.Range("K137").Value = 11474.0001
$ws.Range("L137").Value = 6966894
$ws.Range("M137").Value = -8924.000100000001
$ws.Range("N137").Value = -6971994

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6539180.5
$ws.Range("I32").Value = 6947854.5
$ws.Range("K32").Value = 6947854.5
$ws.Range("M32").Value = -6947567.5
$ws.Range("H63").Value = 23540
$ws.Range("J63").Value = 37269
$ws.Range("L63").Value = 37269
$ws.Range("N63").Value = -38641
$ws.Range("H66").Value = 23540
$ws.Range("J66").Value = 37269
$ws.Range("L66").Value = 186345
$ws.Range("N66").Value = -193209
$ws.Range("H102").Value = 42225.777
$ws.Range("I102").Value = 47254.75
$ws.Range("J102").Value = 1994
$ws.Range("K102").Value = 47254.75
$ws.Range("L102").Value = 1994
$ws.Range("M102").Value = -45632.75
$ws.Range("N102").Value = -5238
$ws.Range("H122").Value = 1918.6666
$ws.Range("I122").Value = 1606
$ws.Range("J122").Value = 2075
$ws.Range("K122").Value = 4818
$ws.Range("L122").Value = 6225
$ws.Range("M122").Value = -2368
$ws.Range("N122").Value = -11125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 9999.5
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 4666.6665
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 4666.6665
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 4666.6665
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = -5136.6665
$ws.Range("H29").Value = 29955
$ws.Range("J29").Value = 29955
$ws.Range("L29").Value = 29955
$ws.Range("N29").Value = -30541
$ws.Range("H31").Value = 113843.42
$ws.Range("J31").Value = 38262.8
$ws.Range("L31").Value = 38262.8
$ws.Range("N31").Value = -38852.8
$ws.Range("H34").Value = 113843.42
$ws.Range("J34").Value = 38262.8
$ws.Range("L34").Value = 38262.8
$ws.Range("N34").Value = -38666.8
$ws.Range("H86").Value = 8198
$ws.Range("I86").Value = 7466.25
$ws.Range("J86").Value = 8929.75
$ws.Range("K86").Value = 7466.25
$ws.Range("L86").Value = 8929.75
$ws.Range("M86").Value = -6343.25
$ws.Range("N86").Value = -11175.75
$ws.Range("H89").Value = 8198
$ws.Range("I89").Value = 7466.25
$ws.Range("J89").Value = 8929.75
$ws.Range("K89").Value = 37331.25
$ws.Range("L89").Value = 44648.75
$ws.Range("M89").Value = -31715.25
$ws.Range("N89").Value = -55880.75
$ws.Range("H134").Value = 7565.32
$ws.Range("I134").Value = 7755.0835
$ws.Range("K134").Value = 23265.2505
$ws.Range("M134").Value = -20730.2505

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 127
$ws.Range("I2").Value = 62.5
$ws.Range("J2").Value = 165.7
$ws.Range("K2").Value = 375
$ws.Range("L2").Value = 994.1999999999999
$ws.Range("M2").Value = -262
$ws.Range("N2").Value = -1220.2
$ws.Range("H22").Value = 933.1667
$ws.Range("I22").Value = 149.75
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 449.25
$ws.Range("L22").Value = 7500
$ws.Range("M22").Value = -280.25
$ws.Range("N22").Value = -7838
$ws.Range("H27").Value = 933.1667
$ws.Range("I27").Value = 149.75
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 449.25
$ws.Range("L27").Value = 7500
$ws.Range("M27").Value = -347.25
$ws.Range("N27").Value = -7704
$ws.Range("H104").Value = 4100
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 3000
$ws.Range("M107").Value = -1080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 4021.4
$ws.Range("I18").Value = 2526.75
$ws.Range("K18").Value = 2526.75
$ws.Range("M18").Value = -2233.75
$ws.Range("H43").Value = 3594.111
$ws.Range("J43").Value = 18000
$ws.Range("L43").Value = 18000
$ws.Range("N43").Value = -18302
$ws.Range("H46").Value = 13973
$ws.Range("I46").Value = 5237.375
$ws.Range("J46").Value = 27950
$ws.Range("K46").Value = 5237.375
$ws.Range("L46").Value = 27950
$ws.Range("M46").Value = -5081.375
$ws.Range("N46").Value = -28262
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").Value = ""
$ws.Range("H80").Value = 323784
$ws.Range("I80").Value = 430252.75
$ws.Range("J80").Value = 4377.75
$ws.Range("K80").Value = 430252.75
$ws.Range("L80").Value = 4377.75
$ws.Range("M80").Value = -429254.75
$ws.Range("N80").Value = -6373.75
$ws.Range("H83").Value = 323784
$ws.Range("I83").Value = 430252.75
$ws.Range("J83").Value = 4377.75
$ws.Range("K83").Value = 2151263.75
$ws.Range("L83").Value = 21888.75
$ws.Range("M83").Value = -2146271.75
$ws.Range("N83").Value = -31872.75
$ws.Range("H126").Value = 1669586.5
$ws.Range("J126").Value = 3737.25
$ws.Range("L126").Value = 11211.75
$ws.Range("N126").Value = -16151.75
$ws.Range("H132").Value = 1516121.6
$ws.Range("I132").Value = 2418835.5
$ws.Range("K132").Value = 7256506.5
$ws.Range("M132").Value = -7253976.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4124.2666
$ws.Range("I7").Value = 3988.6667
$ws.Range("J7").Value = 4666.6665
$ws.Range("K7").Value = 3988.6667
$ws.Range("L7").Value = 4666.6665
$ws.Range("M7").Value = -3876.6667
$ws.Range("N7").Value = -4890.6665
$ws.Range("H17").Value = 683.3333
$ws.Range("I17").Value = 650
$ws.Range("J17").Value = 700
$ws.Range("K17").Value = 650
$ws.Range("L17").Value = 700
$ws.Range("M17").Value = -480
$ws.Range("N17").Value = -1040
$ws.Range("H40").Value = 4689.2354
$ws.Range("I40").Value = 4669.8125
$ws.Range("K40").Value = 4669.8125
$ws.Range("M40").Value = -4533.8125
$ws.Range("H46").Value = 1285.3334
$ws.Range("J46").Value = 2001
$ws.Range("L46").Value = 2001
$ws.Range("N46").Value = -2377
$ws.Range("H61").Value = 3500
$ws.Range("I61").Value = 3250
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 3250
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -3048
$ws.Range("N61").Value = -4404
$ws.Range("H113").Value = 3500
$ws.Range("I113").Value = 3250
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 3250
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -1080
$ws.Range("N113").Value = -8340
$ws.Range("H126").Value = 4124.2666
$ws.Range("I126").Value = 3988.6667
$ws.Range("J126").Value = 4666.6665
$ws.Range("K126").Value = 11966.0001
$ws.Range("L126").Value = 13999.9995
$ws.Range("M126").Value = -9496.000100000001
$ws.Range("N126").Value = -18939.9995
$ws.Range("H132").Value = 4354835
$ws.Range("I132").Value = 17398844
$ws.Range("J132").Value = 6831.3335
$ws.Range("K132").Value = 52196532
$ws.Range("L132").Value = 20494.0005
$ws.Range("M132").Value = -52194002
$ws.Range("N132").Value = -25554.0005
$ws.Range("H138").Value = 88000
$ws.Range("J138").Value = 88000
$ws.Range("L138").Value = 88000
$ws.Range("N138").Value = -98280
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5753642.5
$ws.Range("I132").Value = 6101288
$ws.Range("J132").Value = 17500
$ws.Range("K132").Value = 18303864
$ws.Range("L132").Value = 52500
$ws.Range("M132").Value = -18301334
$ws.Range("N132").Value = -52500
